# Generate Report for Handoff
# Adds two new handed-off files (33ded07e... and 36f66eed...) to the
# localization status report: one new row each on "Overview", "zh-cn"
# and "de-de" sheets, extending their backing tables accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: new rows 4 & 5
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = "33ded07e-c313-4ad1-841d-d92b1e1f4068.md"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("D4").Value = ""
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = "2016-08-13 14:49:33"
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Range("A5").Value = "36f66eed-9d63-451e-97ac-7e7e0964ff41.md"
$wsOverview.Range("C5").Value = ".md"
$wsOverview.Range("D5").Value = ""
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = "2016-08-13 14:49:33"
$wsOverview.Range("G5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/placeholder/e2e/33ded07e-c313-4ad1-841d-d92b1e1f4068.md", "", "", "e2e\33ded07e-c313-4ad1-841d-d92b1e1f4068.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/placeholder/e2e/36f66eed-9d63-451e-97ac-7e7e0964ff41.md", "", "", "e2e\36f66eed-9d63-451e-97ac-7e7e0964ff41.md") | Out-Null

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G5"))

# ---------------------------------------------------------------------
# zh-cn sheet: new rows 4 & 5
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Ready for handoff"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "False"
$wsZhCn.Range("G4").Value = "33ded07e-c313-4ad1-841d-d92b1e1f4068.87d6da3b777a9a6a55dc55760778cb20f41c0d71.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2016-08-13 14:49:25"
$wsZhCn.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I4").Value = ""
$wsZhCn.Range("J4").Value = ""
$wsZhCn.Range("K4").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L4").Value = ""
$wsZhCn.Range("M4").Value = "True"
$wsZhCn.Range("N4").Value = ""
$wsZhCn.Range("O4").Value = "False"
$wsZhCn.Range("P4").Value = ""

$wsZhCn.Range("B5").Value = ".md"
$wsZhCn.Range("C5").Value = "Ready for handoff"
$wsZhCn.Range("D5").Value = "e2e"
$wsZhCn.Range("E5").Value = "ht"
$wsZhCn.Range("F5").Value = "False"
$wsZhCn.Range("G5").Value = "36f66eed-9d63-451e-97ac-7e7e0964ff41.c59bcfa99ce50b02421fb9afb96e73b95fe449b2.zh-cn.xlf"
$wsZhCn.Range("H5").Value = "2016-08-13 14:49:25"
$wsZhCn.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I5").Value = ""
$wsZhCn.Range("J5").Value = ""
$wsZhCn.Range("K5").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L5").Value = ""
$wsZhCn.Range("M5").Value = "True"
$wsZhCn.Range("N5").Value = ""
$wsZhCn.Range("O5").Value = "False"
$wsZhCn.Range("P5").Value = ""

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/placeholder/e2e/33ded07e-c313-4ad1-841d-d92b1e1f4068.md", "", "", "33ded07e-c313-4ad1-841d-d92b1e1f4068.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/placeholder/e2e/36f66eed-9d63-451e-97ac-7e7e0964ff41.md", "", "", "36f66eed-9d63-451e-97ac-7e7e0964ff41.md") | Out-Null

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P5"))

# ---------------------------------------------------------------------
# de-de sheet: new rows 4 & 5
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Ready for handoff"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "False"
$wsDeDe.Range("G4").Value = "33ded07e-c313-4ad1-841d-d92b1e1f4068.87d6da3b777a9a6a55dc55760778cb20f41c0d71.de-de.xlf"
$wsDeDe.Range("H4").Value = "2016-08-13 14:49:33"
$wsDeDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I4").Value = ""
$wsDeDe.Range("J4").Value = ""
$wsDeDe.Range("K4").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L4").Value = ""
$wsDeDe.Range("M4").Value = "True"
$wsDeDe.Range("N4").Value = ""
$wsDeDe.Range("O4").Value = "False"
$wsDeDe.Range("P4").Value = ""

$wsDeDe.Range("B5").Value = ".md"
$wsDeDe.Range("C5").Value = "Ready for handoff"
$wsDeDe.Range("D5").Value = "e2e"
$wsDeDe.Range("E5").Value = "ht"
$wsDeDe.Range("F5").Value = "False"
$wsDeDe.Range("G5").Value = "36f66eed-9d63-451e-97ac-7e7e0964ff41.c59bcfa99ce50b02421fb9afb96e73b95fe449b2.de-de.xlf"
$wsDeDe.Range("H5").Value = "2016-08-13 14:49:33"
$wsDeDe.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I5").Value = ""
$wsDeDe.Range("J5").Value = ""
$wsDeDe.Range("K5").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L5").Value = ""
$wsDeDe.Range("M5").Value = "True"
$wsDeDe.Range("N5").Value = ""
$wsDeDe.Range("O5").Value = "False"
$wsDeDe.Range("P5").Value = ""

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/placeholder/e2e/33ded07e-c313-4ad1-841d-d92b1e1f4068.md", "", "", "33ded07e-c313-4ad1-841d-d92b1e1f4068.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/placeholder/e2e/36f66eed-9d63-451e-97ac-7e7e0964ff41.md", "", "", "36f66eed-9d63-451e-97ac-7e7e0964ff41.md") | Out-Null

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P5"))

Write-Output "edit complete"
